$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new row 14 with data, mirroring the existing rows' layout (A:N)
$ws.Range("A14").Value = 42622.891446759262

$ws.Range("B14").Value = 14
$ws.Range("C14").Value = 58
$ws.Range("D14").Value = 41
$ws.Range("E14").Value = 58
$ws.Range("F14").Value = 100
$ws.Range("G14").Value = 19787
$ws.Range("H14").Value = 10163
$ws.Range("I14").Value = 1498
$ws.Range("J14").Value = 264
$ws.Range("K14").Value = 189
$ws.Range("L14").Value = 0
$ws.Range("M14").Value = 5
$ws.Range("N14").Value = "Noun"
